$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-8 (columns D, I, J, K, L, M, N, P, Q)
# All other columns (A,B,C,E,F,G,H,O,R) remain unchanged.

$ws.Range("D2").Value = 44235
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("N2").Value = "$/bandeja 18 kilos"
$ws.Range("P2").Value = 778
$ws.Range("Q2").Value = 18

$ws.Range("D3").Value = 44235
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("N3").Value = "$/bandeja 18 kilos"
$ws.Range("P3").Value = 667
$ws.Range("Q3").Value = 18

$ws.Range("D4").Value = 44235
$ws.Range("I4").Value = "Tercera"
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("N4").Value = "$/bandeja 18 kilos"
$ws.Range("P4").Value = 556
$ws.Range("Q4").Value = 18

$ws.Range("D5").Value = 44424
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 75
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 18000
$ws.Range("N5").Value = "$/caja 15 kilos"
$ws.Range("P5").Value = 1200
$ws.Range("Q5").Value = 15

$ws.Range("D6").Value = 44424
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 12000
$ws.Range("N6").Value = "$/caja 15 kilos"
$ws.Range("P6").Value = 800
$ws.Range("Q6").Value = 15

$ws.Range("D7").Value = 44242
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("N7").Value = "$/bandeja 18 kilos"
$ws.Range("P7").Value = 722
$ws.Range("Q7").Value = 18

$ws.Range("D8").Value = 44242
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 10000
$ws.Range("N8").Value = "$/bandeja 18 kilos"
$ws.Range("P8").Value = 556
$ws.Range("Q8").Value = 18
